# Adds rows 11-13 to the "Artfynd" worksheet (new species observations),
# matching the source export format where every field is written as text
# unless it is a genuine number/boolean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 11 ----
$ws.Range("A11").Value = 112244426
$ws.Range("B11").Value = 77550
$ws.Range("C11").Value = 'Ovaliderad'
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 185
$ws.Range("F11").Value = 'Violettgrå tagellav'
$ws.Range("G11").Value = 'Bryoria nadvornikiana'
$ws.Range("H11").Value = '(Gyeln.) Brodo & D.Hawksw.'
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = ""
$ws.Range("I11").Style = "Normal"
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = ""
$ws.Range("K11").Style = "Normal"
$ws.Range("P11").Value = 'Berg-Andersberget, Dlr'
$ws.Range("Q11").Value = 523094
$ws.Range("R11").Value = 6739613
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 'Dalarna'
$ws.Range("U11").Value = 'Falun'
$ws.Range("V11").Value = 'Dalarna'
$ws.Range("W11").Value = 'Bjursås'
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = '2023-09-22'
$ws.Range("Y11").Style = "Normal"
$ws.Range("Z11").NumberFormat = "@"
$ws.Range("Z11").Value = '11:11'
$ws.Range("Z11").Style = "Normal"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = '2023-09-22'
$ws.Range("AA11").Style = "Normal"
$ws.Range("AB11").NumberFormat = "@"
$ws.Range("AB11").Value = '11:11'
$ws.Range("AB11").Style = "Normal"
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AJ11").Value = 'gran'
$ws.Range("AK11").Value = 'Picea abies'
$ws.Range("AO11").Value = 'Picea abies'
$ws.Range("AT11").NumberFormat = "@"
$ws.Range("AT11").Value = ""
$ws.Range("AT11").Style = "Normal"
$ws.Range("AW11").Value = 'Uno Skog'
$ws.Range("AX11").Value = 'Uno Skog, Anton Björk'
$ws.Range("AY11").NumberFormat = "@"
$ws.Range("AY11").Value = ""
$ws.Range("AY11").Style = "Normal"

# ---- Row 12 ----
$ws.Range("A12").Value = 112243622
$ws.Range("B12").Value = 89745
$ws.Range("C12").Value = 'Ovaliderad'
$ws.Range("D12").Value = 'VU'
$ws.Range("E12").Value = 2062
$ws.Range("F12").Value = 'Ulltickeporing'
$ws.Range("G12").Value = 'Skeletocutis brevispora'
$ws.Range("H12").Value = 'Niemelä'
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = ""
$ws.Range("I12").Style = "Normal"
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = ""
$ws.Range("K12").Style = "Normal"
$ws.Range("P12").Value = 'Berg-Andersberget, Dlr'
$ws.Range("Q12").Value = 523006
$ws.Range("R12").Value = 6739484
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = 'Dalarna'
$ws.Range("U12").Value = 'Falun'
$ws.Range("V12").Value = 'Dalarna'
$ws.Range("W12").Value = 'Bjursås'
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = '2023-09-22'
$ws.Range("Y12").Style = "Normal"
$ws.Range("Z12").NumberFormat = "@"
$ws.Range("Z12").Value = '10:21'
$ws.Range("Z12").Style = "Normal"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = '2023-09-22'
$ws.Range("AA12").Style = "Normal"
$ws.Range("AB12").NumberFormat = "@"
$ws.Range("AB12").Value = '10:21'
$ws.Range("AB12").Style = "Normal"
$ws.Range("AC12").Value = 'På ytmurken granlåga med delvis avfallande bark och insektsgnag i veden från tiden då granen stod upp.'
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AJ12").Value = 'ullticka'
$ws.Range("AK12").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("AO12").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("AT12").NumberFormat = "@"
$ws.Range("AT12").Value = ""
$ws.Range("AT12").Style = "Normal"
$ws.Range("AW12").Value = 'Uno Skog'
$ws.Range("AX12").Value = 'Uno Skog, Anton Björk'
$ws.Range("AY12").NumberFormat = "@"
$ws.Range("AY12").Value = ""
$ws.Range("AY12").Style = "Normal"

# ---- Row 13 ----
$ws.Range("A13").Value = 112243230
$ws.Range("B13").Value = 89405
$ws.Range("C13").Value = 'Ovaliderad'
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = 'Ullticka'
$ws.Range("G13").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H13").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = ""
$ws.Range("I13").Style = "Normal"
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value = ""
$ws.Range("K13").Style = "Normal"
$ws.Range("P13").Value = 'Berg-Andersberget, Dlr'
$ws.Range("Q13").Value = 523006
$ws.Range("R13").Value = 6739484
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = 'Dalarna'
$ws.Range("U13").Value = 'Falun'
$ws.Range("V13").Value = 'Dalarna'
$ws.Range("W13").Value = 'Bjursås'
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = '2023-09-22'
$ws.Range("Y13").Style = "Normal"
$ws.Range("Z13").NumberFormat = "@"
$ws.Range("Z13").Value = '10:16'
$ws.Range("Z13").Style = "Normal"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = '2023-09-22'
$ws.Range("AA13").Style = "Normal"
$ws.Range("AB13").NumberFormat = "@"
$ws.Range("AB13").Value = '10:16'
$ws.Range("AB13").Style = "Normal"
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AJ13").Value = 'gran'
$ws.Range("AK13").Value = 'Picea abies'
$ws.Range("AO13").Value = 'Picea abies'
$ws.Range("AT13").NumberFormat = "@"
$ws.Range("AT13").Value = ""
$ws.Range("AT13").Style = "Normal"
$ws.Range("AW13").Value = 'Uno Skog'
$ws.Range("AX13").Value = 'Uno Skog, Anton Björk'
$ws.Range("AY13").NumberFormat = "@"
$ws.Range("AY13").Value = ""
$ws.Range("AY13").Style = "Normal"
